# Updated symbol list on Thu Jan  5 20:32:18 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) figures (and, for the rows that lost the
# "HuobiToken-rank" entry, shift Coin/Link/Price/Volume up one row) to
# match the latest coinranking.com scrape.
#
# Values that look numeric (plain numbers or "NN.NN%") are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inlineStr/text cells) instead of auto-converting them to
# numbers or percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''257.64'
$ws.Range("E2").Value = '''0.05%'
$ws.Range("D3").Value = '''27.24'
$ws.Range("E3").Value = '''-1.38%'
$ws.Range("E4").Value = '''-12.00%'
$ws.Range("D5").Value = '''0.05887'
$ws.Range("E5").Value = '''-0.61%'
$ws.Range("D6").Value = '''6.645'
$ws.Range("E6").Value = '''-0.71%'
$ws.Range("D7").Value = '''0.8580'
$ws.Range("E7").Value = '''-1.12%'
$ws.Range("D8").Value = '''0.9424'
$ws.Range("E8").Value = '''-8.38%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1410'
$ws.Range("E9").Value = '''-0.70%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.03983'
$ws.Range("E10").Value = '''10.59%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07089'
$ws.Range("E11").Value = '''-1.48%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03177'
$ws.Range("E12").Value = '''-2.70%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09157'
$ws.Range("E13").Value = '''-0.64%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001536'
$ws.Range("E14").Value = '''-1.01%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '''0.0006029'
$ws.Range("E15").Value = '''-0.72%'
$ws.Range("D16").Value = '''0.006006'
$ws.Range("E16").Value = '''2.10%'
$ws.Range("E17").Value = '''0.62%'
$ws.Range("D18").Value = '''3.207'
$ws.Range("E18").Value = '''-1.98%'
$ws.Range("E19").Value = '''0.89%'
$ws.Range("E21").Value = '''-1.17%'
$ws.Range("D22").Value = '''3.907'
$ws.Range("E22").Value = '''10.86%'
$ws.Range("D23").Value = '''0.04226'
$ws.Range("E23").Value = '''1.56%'
$ws.Range("D24").Value = '''0.001221'
$ws.Range("E24").Value = '''0.23%'
$ws.Range("D25").Value = '''0.004293'
$ws.Range("E25").Value = '''-5.07%'
$ws.Range("D26").Value = '''0.0001200'
$ws.Range("E26").Value = '''-0.08%'
$ws.Range("D27").Value = '''0.0001937'
$ws.Range("E27").Value = '''-0.12%'
$ws.Range("D40").Value = '''0.03853'
$ws.Range("E40").Value = '''0.76%'
$ws.Range("D41").Value = '''0.006218'
$ws.Range("E41").Value = '''14.88%'
$ws.Range("D42").Value = '''0.1104'
$ws.Range("E42").Value = '''-0.05%'
$ws.Range("D43").Value = '''0.002409'
$ws.Range("E43").Value = '''26.74%'
$ws.Range("D44").Value = '''0.01176'
$ws.Range("E44").Value = '''19.69%'
$ws.Range("D45").Value = '''0.00005457'
$ws.Range("E45").Value = '''0.44%'
$ws.Range("E46").Value = '''-0.06%'
$ws.Range("D47").Value = '''0.05999'
$ws.Range("E47").Value = '''-45.04%'
$ws.Range("D48").Value = '''0.1281'
$ws.Range("E48").Value = '''5,825.04%'
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("E49").Value = '''-0.06%'
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("E50").Value = '''-0.06%'
